$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 updates
$ws.Range("G11").Value = 4.1
$ws.Range("H11").Value = 3.4
$ws.Range("I11").Value = 1.85
$ws.Range("J11").Value = 4.75
$ws.Range("K11").Value = 2.05
$ws.Range("L11").Value = 2.6
$ws.Range("M11").Value = 1.07
$ws.Range("N11").Value = 9
$ws.Range("O11").Value = 1.4
$ws.Range("P11").Value = 2.75
$ws.Range("Q11").Value = 2.2
$ws.Range("R11").Value = 1.65
$ws.Range("S11").Value = 1.44
$ws.Range("T11").Value = 2.63
$ws.Range("U11").Value = 2
$ws.Range("V11").Value = 1.73
$ws.Range("W11").Value = 10
$ws.Range("Z11").Value = 41
$ws.Range("AB11").Value = 41
$ws.Range("AC11").Value = 8
$ws.Range("AH11").Value = 6
$ws.Range("AJ11").Value = 9
$ws.Range("AL11").Value = 17
$ws.Range("AO11").Value = 23
$ws.Range("AP11").Value = 34
$ws.Range("AQ11").Value = 81
$ws.Range("AR11").Value = 126
$ws.Range("AS11").Value = 301
$ws.Range("AT11").Value = 2.63
$ws.Range("AU11").Value = 8.5
$ws.Range("AX11").Value = 10
$ws.Range("AY11").Value = 23
$ws.Range("AZ11").Value = 34
$ws.Range("BA11").Value = 51
$ws.Range("BB11").Value = 201

# Row 12 updates
$ws.Range("Q12").Value = 2.1
$ws.Range("R12").Value = 1.7
